$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# C2: was stored as a raw number (10000092508), change it to a text value
# matching the rest of the reg_no column (shared string, not numeric).
$ws.Range("C2").Value = "10000092508"

# Y2:Y30 (shift2 column): change shift value from "I" to "II"
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 25).Value = "II"
}

# Update the view/selection to match: active cell Y2, selected range Y2:Y30
$ws.Range("Y2:Y30").Select()
